# power meter remote 외형 구현 offset 구현중
# Rebuild the "Select Freq" unit (GHz -> MHz) and the M-column offset ladder
# on the "config" sheet, extending the used range from row 13 to row 104.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("config")

# --- M2: unit label GHz -> MHz -----------------------------------------
$ws.Range("M2").Value = "MHz"

# --- Row 3: per-channel numbers change ----------------------------------
$ws.Range("H3").Value = 18
$ws.Range("I3").Value = 20
$ws.Range("M3").Value = 700
$ws.Range("N3").Value = 10

# --- Row 4: starting offset ---------------------------------------------
$ws.Range("M4").Value = 1000

# --- Row 5: step value (N5) + first running-offset formula --------------
$ws.Range("N5").Value = 50
$ws.Range("M5").Formula = '=M4+$N$5'

# --- M6:M69 shared formula group: each = previous + step ----------------
$ws.Range("M6:M69").Formula = '=M5+$N$5'

# --- M70:M122 second shared formula group (kept even though the sheet is
#     later trimmed back to row 104, matching the authored workbook) -----
$ws.Range("M70:M122").Formula = '=M69+$N$5'

# The working range was trimmed back down to row 104 after being built out
# further, leaving dimension/used-range at A1:R104.
$ws.Rows("105:122").Delete()

# --- Selection ends up on L91 after scrolling down to review the ladder --
[void]$ws.Range("L91").Select()
